# Generate Report for Handback
#
# Row 6 on the "zh-cn" and "de-de" worksheets corresponds to file
# 0df7930f-9db3-4100-9ef1-e043c9694049.md. A new handback (translated
# target file) has now come in for that row, but it is *stale* relative
# to the latest handoff, so we record the target file, the handback
# xliff file, the handback datetime, and a validation error message, and
# widen the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$targetFile    = "0df7930f-9db3-4100-9ef1-e043c9694049.md"
$currentCommit = "e6ab0b68623f6dc00106ea9c12aeddea169affd0"
$latestCommit  = "b71a1fc9633a25fb6ea8166a67694a51878189b1"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$currentCommit/e2e/$targetFile, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$latestCommit/e2e/$targetFile."

# Column "P" (Error Detail) needs to be wide enough to show the message.
# (39.1666... set here renders as stored width 40, matching the other
# autosized-to-40 columns on this sheet.)
$wideColumnWidth = 39.16666666666667

$sheets = @(
    @{ Name = "zh-cn"; Repo = "ol-test0-zhcn"; HandbackFile = "0df7930f-9db3-4100-9ef1-e043c9694049.c4507590df30fa4b347fd4286b12c2ed353f53d2.zh-cn.xlf"; HandbackDateTime = "2016-08-27 12:40:54" },
    @{ Name = "de-de"; Repo = "ol-test0-dede"; HandbackFile = "0df7930f-9db3-4100-9ef1-e043c9694049.c4507590df30fa4b347fd4286b12c2ed353f53d2.de-de.xlf"; HandbackDateTime = "2016-08-27 12:41:02" }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen the Error Detail column (16th column = "P").
    $ws.Columns.Item(16).ColumnWidth = $wideColumnWidth

    # I6 - Latest Target File: link to the handed-back file in the
    # language-specific repo.
    $targetUrl = "https://github.com/OpenLocalizationTestOrg/$($info.Repo)/blob/$currentCommit/e2e/$targetFile"
    $ws.Hyperlinks.Add($ws.Range("I6"), $targetUrl, "", "", $targetFile)

    # J6 - Latest Handback File.
    $ws.Range("J6").Value = $info.HandbackFile

    # K6 - Latest Handback DateTime.
    $ws.Range("K6").Value = $info.HandbackDateTime

    # P6 - Error Detail.
    $ws.Range("P6").Value = $errorDetail
}

Write-Host "Handback report generated for row 6 on zh-cn and de-de sheets."
